$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Insert four new blank columns before the existing ExpPoints column (C).
# This shifts the old column C (ExpPoints, with its data + header) to
# column G, and the newly inserted C:F columns inherit column C's
# original formatting (including the bold/bordered header style) for
# row 1.
$ws.Range("C1:F1").EntireColumn.Insert()

# Re-label the new header cells; ExpPoints (now in G1) already carries
# over its original text and style from the insert, so it needs no
# further changes.
$ws.Cells.Item(1, 3).Value = "WIN"
$ws.Cells.Item(1, 4).Value = "TOP2"
$ws.Cells.Item(1, 5).Value = "TOP4"
$ws.Cells.Item(1, 6).Value = "RELEGATION"

# The data rows for the 4 new columns are placeholders for now (to be
# filled in later by the Monte Carlo simulation) - touch each cell so it
# is persisted as an (empty) cell in the sheet, without picking up any
# new formatting/style.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 3; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Style = "Normal"
    }
}
